$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (rows 231-233) to append below the existing 230 rows.
$rows = @(
    @{
        A = 6
        B = "Mercado Mayorista Lo Valledor de Santiago"
        C = "Metropolitana"
        D = 44911
        E = 13
        F = "Fruta"
        G = 100101
        H = "Berries"
        I = 100101004
        J = "Frambuesa"
        K = "Sin especificar"
        L = "Especial"
        M = 200
        N = 7000
        O = 7000
        P = 7000
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Curicó"
        S = 3500
        T = 2
    },
    @{
        A = 6
        B = "Mercado Mayorista Lo Valledor de Santiago"
        C = "Metropolitana"
        D = 44911
        E = 13
        F = "Fruta"
        G = 100101
        H = "Berries"
        I = 100101004
        J = "Frambuesa"
        K = "Sin especificar"
        L = "Especial"
        M = 350
        N = 7000
        O = 7000
        P = 7000
        Q = "`$/bandeja 2 kilos"
        R = "Región del Maule"
        S = 3500
        T = 2
    },
    @{
        A = 6
        B = "Mercado Mayorista Lo Valledor de Santiago"
        C = "Metropolitana"
        D = 44911
        E = 13
        F = "Fruta"
        G = 100101
        H = "Berries"
        I = 100101004
        J = "Frambuesa"
        K = "Sin especificar"
        L = "Primera"
        M = 250
        N = 6000
        O = 6000
        P = 6000
        Q = "`$/bandeja 2 kilos"
        R = "Provincia de Curicó"
        S = 3000
        T = 2
    }
)

$startRow = 231
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    foreach ($col in $cols) {
        $cell = $ws.Range($col + $r)
        $cell.Value = $data[$col]
        if ($col -eq "D") {
            $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
        }
    }
}
